$d = $word.ActiveDocument
$found = $d.Content.Find.Execute(". Typy testów", $true, $false, $false, $false, $false,
                         $true, 1, $false, " Typy testów", 2)
Write-Output "Find result: $found"
